# Update the division-problem answers in the first (and only) table.
# Each data row (1, 5, 9, 13, 17) holds five worked-out division
# problems, one per column. We address cells by (row, column) rather
# than Find/Replace so that values which happen to collide between the
# "old" and "new" sets (e.g. "687÷8=85, 7" is both replaced in one cell
# and introduced in another) can never cross-contaminate each other.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="985÷3=328, 1"},
    @{Row=1;  Col=2; Text="632÷8=79, 0"},
    @{Row=1;  Col=3; Text="874÷5=174, 4"},
    @{Row=1;  Col=4; Text="570÷7=81, 3"},
    @{Row=1;  Col=5; Text="945÷8=118, 1"},

    @{Row=5;  Col=1; Text="732÷2=366, 0"},
    @{Row=5;  Col=2; Text="324÷7=46, 2"},
    @{Row=5;  Col=3; Text="494÷5=98, 4"},
    @{Row=5;  Col=4; Text="873÷8=109, 1"},
    @{Row=5;  Col=5; Text="832÷3=277, 1"},

    @{Row=9;  Col=1; Text="632÷4=158, 0"},
    @{Row=9;  Col=2; Text="311÷5=62, 1"},
    @{Row=9;  Col=3; Text="810÷5=162, 0"},
    @{Row=9;  Col=4; Text="565÷9=62, 7"},
    @{Row=9;  Col=5; Text="786÷3=262, 0"},

    @{Row=13; Col=1; Text="705÷8=88, 1"},
    @{Row=13; Col=2; Text="687÷8=85, 7"},
    @{Row=13; Col=3; Text="875÷6=145, 5"},
    @{Row=13; Col=4; Text="477÷6=79, 3"},
    @{Row=13; Col=5; Text="833÷3=277, 2"},

    @{Row=17; Col=1; Text="196÷2=98, 0"},
    @{Row=17; Col=2; Text="541÷5=108, 1"},
    @{Row=17; Col=3; Text="803÷2=401, 1"},
    @{Row=17; Col=4; Text="207÷5=41, 2"},
    @{Row=17; Col=5; Text="370÷9=41, 1"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
